# Fruta / hortaliza, semanal
# Insert two new weekly observation rows before the existing row 239,
# pushing the existing rows 239:249 down to 241:251, then populate the
# two freshly-inserted rows (239:240) with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above row 239 (existing data shifts down to 241:251)
$ws.Rows.Item(239).Resize(2).Insert()

# New row 239 data
$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44509
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = 100112032
$ws.Cells.Item(239, 7).Value = "Zapallo italiano"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 80
$ws.Cells.Item(239, 11).Value = 4000
$ws.Cells.Item(239, 12).Value = 4000
$ws.Cells.Item(239, 13).Value = 4000
$ws.Cells.Item(239, 14).Value = "$/caja 36 unidades"
$ws.Cells.Item(239, 15).Value = "Limache"
$ws.Cells.Item(239, 16).Value = 111
$ws.Cells.Item(239, 17).Value = 36
$ws.Cells.Item(239, 18).Value = "Hortaliza"

# New row 240 data
$ws.Cells.Item(240, 1).Value = 3
$ws.Cells.Item(240, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(240, 3).Value = "Coquimbo"
$ws.Cells.Item(240, 4).Value = 44509
$ws.Cells.Item(240, 5).Value = 5
$ws.Cells.Item(240, 6).Value = 100112032
$ws.Cells.Item(240, 7).Value = "Zapallo italiano"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 110
$ws.Cells.Item(240, 11).Value = 8000
$ws.Cells.Item(240, 12).Value = 8500
$ws.Cells.Item(240, 13).Value = 8227
$ws.Cells.Item(240, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(240, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(240, 16).Value = 118
$ws.Cells.Item(240, 17).Value = 70
$ws.Cells.Item(240, 18).Value = "Hortaliza"
